# Refresh the cryptos price/volume table (columns D = Price, E = Volume(1h))
# for rows 2-51. Values are written as literal text (a leading "'" is used
# for numeric-looking strings, matching how Excel stores an apostrophe-
# prefixed entry) so formatting such as trailing zeros ("1.000",
# "0.07740") and the "000...0" style prices is preserved exactly, instead
# of being auto-converted to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.394.27"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.847.30"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'240.47"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'0.6286"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07605"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "'0.2928"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").Value = "'24.48"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'0.07740"
$ws.Range("D12").Value = "1.848.68"
$ws.Range("E12").Value = "  -6.87%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'0.00001082"
$ws.Range("E14").Value = "  +9.07%  "
$ws.Range("D15").Value = "'0.6788"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "'83.74"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "2.097.43"
$ws.Range("E17").Value = "  -7.37%  "
$ws.Range("D18").Value = "'6.175"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "29.411.19"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'228.50"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'7.472"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'157.29"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "'0.1394"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "'8.352"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "'1.465"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'1.297"
$ws.Range("E30").Value = "  +3.69%  "
$ws.Range("D31").Value = "'0.05580"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("D32").Value = "'4.099"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "'4.030"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'1.846"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'1.155"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'0.7092"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "'2.585"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "1.233.07"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").Value = "'2.773"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "'0.01800"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "'6.410"
$ws.Range("E41").Value = "  +5.21%  "
$ws.Range("D42").Value = "'0.9058"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'101.78"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'66.03"
$ws.Range("D46").Value = "'0.00000000121"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").Value = "'7.173"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "'0.4018"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "'9.020"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("D50").Value = "'1.676"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").Value = "'0.1121"
$ws.Range("E51").Value = "  -0.46%  "
